$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
